# Auto-generated edit script: refresh crypto price/volume data
# and reorder a few rows, per commit "Updated cryptos list on Mon Jul 15 12:14:15 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.624.12'
$ws.Range('E2').Value = '  +4.50%  '
$ws.Range('D3').Value = '3.343.22'
$ws.Range('E3').Value = '  +4.59%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'564.91"
$ws.Range('E5').Value = '  +5.45%  '
$ws.Range('D6').Value = "'152.37"
$ws.Range('E6').Value = '  +5.43%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').Value = '3.346.87'
$ws.Range('E8').Value = '  +4.49%  '
$ws.Range('D9').Value = "'0.534"
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('D10').Value = "'7.42"
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  +4.59%  '
$ws.Range('D12').Value = "'0.437"
$ws.Range('E12').Value = '  +1.96%  '
$ws.Range('D13').Value = '3.924.45'
$ws.Range('E13').Value = '  +4.60%  '
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').Value = "'26.83"
$ws.Range('E15').Value = '  +4.10%  '
$ws.Range('D16').Value = "'0.0000179"
$ws.Range('E16').Value = '  +3.76%  '
$ws.Range('D17').Value = '62.631.03'
$ws.Range('E17').Value = '  +4.36%  '
$ws.Range('D18').Value = '3.339.91'
$ws.Range('E18').Value = '  +4.86%  '
$ws.Range('D19').Value = "'6.35"
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = "'13.82"
$ws.Range('E20').Value = '  +5.84%  '
$ws.Range('D21').Value = "'8.41"
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('D22').Value = "'387.62"
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('E24').Value = '  +1.97%  '
$ws.Range('D25').Value = "'69.92"
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = "'0.180"
$ws.Range('E26').Value = '  +5.80%  '
$ws.Range('D27').Value = "'9.14"
$ws.Range('E27').Value = '  +3.92%  '
$ws.Range('E28').Value = '  +6.72%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = "'1.99"
$ws.Range('E30').Value = '  +4.78%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = "'6.47"
$ws.Range('E31').Value = '  +5.97%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'22.94"
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = "'5.58"
$ws.Range('E33').Value = '  +4.35%  '
$ws.Range('E34').Value = '  +9.38%  '
$ws.Range('E35').Value = '  +1.60%  '
$ws.Range('D36').Value = "'1.48"
$ws.Range('E36').Value = '  +10.21%  '
$ws.Range('D37').Value = "'158.88"
$ws.Range('E37').Value = '  +1.40%  '
$ws.Range('E38').Value = '  +12.14%  '
$ws.Range('D39').Value = "'27.04"
$ws.Range('E39').Value = '  +5.93%  '
$ws.Range('E40').Value = '  +5.15%  '
$ws.Range('D41').Value = '2.791.59'
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').Value = "'0.0320"
$ws.Range('E42').Value = '  +8.18%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = "'40.55"
$ws.Range('E43').Value = '  +1.97%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = "'0.744"
$ws.Range('E44').Value = '  +4.04%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = "'4.25"
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('E46').Value = '  +5.06%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = "'22.07"
$ws.Range('E47').Value = '  +7.62%  '
$ws.Range('B48').Value = 'RenzoRestakedETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D48').Value = '3.388.23'
$ws.Range('E48').Value = '  +4.51%  '
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('D50').Value = "'6.31"
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('D51').Value = "'288.41"
$ws.Range('E51').Value = '  +8.74%  '
